$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "'13"
$ws.Range("B2").Value = 5809570196
$ws.Range("C2").Value = "2021-11-21 23:33:21"
$ws.Range("D2").Value = "桃树下的孩子"
$ws.Range("E2").Value = "感谢字幕！第六场的普通版剪辑版双声道的少年/少女记忆虽然别有一番风味，而且剪辑在一起能明显看出很多动作都是同步或者对称的，对照着看特别有感觉，但是能看到清晰独立版本的活着的只有我（？）真的好棒！以及！！！真诚安利大家关注这场犬彦和宫比的互动，包括狼欒神社solo的时候模拟打鼓啦，井户曲摸头啦，魔神曲犬彦拉住宫比的手揽住他的腰等等，互动又多又甜（当然其他场也很甜！我记得有一场，忘记是哪场了狼欒solo时宫比在后面给犬彦比心心！），这对青梅竹马（大概是吧）真的超级超级好嗑！！！"
$ws.Range("F2").Value = 5
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 0

$ws.Range("A3").Value = "'13.1"
$ws.Range("B3").Value = 5809780005
$ws.Range("C3").Value = "2021-11-22 00:05:05"
$ws.Range("D3").Value = "Ponster_"
$ws.Range("E3").Value = "弹幕中的翻译佬！！感谢指正[脱单doge]"
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 5809570196

$ws.Range("A4").Value = "'14"
$ws.Range("B4").Value = 5809565870
$ws.Range("C4").Value = "2021-11-21 23:33:05"
$ws.Range("D4").Value = "我心向云月"
$ws.Range("E4").Value = "西装男跟夜姬一起演太可怕了，jk快远离疯批男[冷]"
$ws.Range("F4").Value = 4
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0

$ws.Range("A5").Value = "'1"
$ws.Range("B5").Value = 5817605996
$ws.Range("C5").Value = "2021-11-23 12:32:48"
$ws.Range("D5").Value = "总攻祁墨宸大人"
$ws.Range("E5").Value = "飞速闻讯而来"
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0

$ws.Range("A6").Value = "'2"
$ws.Range("B6").Value = 5815985159
$ws.Range("C6").Value = "2021-11-23 01:26:12"
$ws.Range("D6").Value = "木容秀吉"
$ws.Range("E6").Value = "这个安可我直呼好家伙"
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0

$ws.Range("A7").Value = "'19"
$ws.Range("B7").Value = 5807519793
$ws.Range("C7").Value = "2021-11-21 18:25:36"
$ws.Range("D7").Value = "墨弦青风"
$ws.Range("E7").Value = "感谢up，up辛苦了[热词系列_吹爆]"
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 0

$ws.Range("A8").Value = "'19.1"
$ws.Range("B8").Value = 5810526423
$ws.Range("C8").Value = "2021-11-22 06:24:28"
$ws.Range("D8").Value = "Ponster_"
$ws.Range("E8").Value = "感谢(=・ω・=)"
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 5807519793

$ws.Range("A9").Value = "'18"
$ws.Range("B9").Value = 5808208456
$ws.Range("C9").Value = "2021-11-21 20:10:25"
$ws.Range("D9").Value = "Ponster_"
$ws.Range("E9").Value = "前面传错了版本，已更正。`n曲目信息、想说的话都在视频里了。`n第一场重制版也已上传，链接https://www.bilibili.com/video/BV1EU4y1u7HA"
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0

$ws.Range("A10").Value = "'17"
$ws.Range("B10").Value = 5808292895
$ws.Range("C10").Value = "2021-11-21 20:23:50"
$ws.Range("D10").Value = "烟云z"
$ws.Range("E10").Value = "太顶了老哥[BW2020_棒棒哦]"
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0

$ws.Range("A11").Value = "'16"
$ws.Range("B11").Value = 5808581488
$ws.Range("C11").Value = "2021-11-21 21:09:35"
$ws.Range("D11").Value = "羽蛇的尾巴尖"
$ws.Range("E11").Value = "迅速缓存"
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 1
$ws.Range("H11").Value = 0

$ws.Range("A12").Value = "'16.1"
$ws.Range("B12").Value = 5810522401
$ws.Range("C12").Value = "2021-11-22 06:24:00"
$ws.Range("D12").Value = "Ponster_"
$ws.Range("E12").Value = "[tv_点赞]"
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 5808581488

$ws.Range("A13").Value = "'15"
$ws.Range("B13").Value = 5809188522
$ws.Range("C13").Value = "2021-11-21 22:40:13"
$ws.Range("D13").Value = "Panic-"
$ws.Range("E13").Value = "草西装男好可怕"
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0

$ws.Range("A14").Value = "'12"
$ws.Range("B14").Value = 5809720823
$ws.Range("C14").Value = "2021-11-21 23:56:28"
$ws.Range("D14").Value = "VirginMary"
$ws.Range("E14").Value = "好耶 来了来了"
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0

$ws.Range("A15").Value = "'11"
$ws.Range("B15").Value = 5809727182
$ws.Range("C15").Value = "2021-11-21 23:58:18"
$ws.Range("D15").Value = "夏空凛冬至"
$ws.Range("E15").Value = "谢谢up！不知道还有没有其他的！"
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = 0

$ws.Range("A16").Value = "'11.1"
$ws.Range("B16").Value = 5809752664
$ws.Range("C16").Value = "2021-11-22 00:01:40"
$ws.Range("D16").Value = "Ponster_"
$ws.Range("E16").Value = "暂时...不会做了吧，原因看视频开头"
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 5809727182

$ws.Range("A17").Value = "'11.2"
$ws.Range("B17").Value = 5809822355
$ws.Range("C17").Value = "2021-11-22 00:13:14"
$ws.Range("D17").Value = "我心向云月"
$ws.Range("E17").Value = "回复 @Ponster_ :可以求其他安可的生肉吗[大哭]"
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 5809727182

$ws.Range("A18").Value = "'11.3"
$ws.Range("B18").Value = 5810008526
$ws.Range("C18").Value = "2021-11-22 00:49:30"
$ws.Range("D18").Value = "Hexachlorocyclohexane"
$ws.Range("E18").Value = "回复 @Ponster_ :同求其他安可的生肉[大哭]"
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 5809727182

$ws.Range("A19").Value = "'10"
$ws.Range("B19").Value = 5810353355
$ws.Range("C19").Value = "2021-11-22 03:00:50"
$ws.Range("D19").Value = "召唤魔术"
$ws.Range("E19").Value = "求上传其他场次生肉[tv_大佬]"
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 0

$ws.Range("A20").Value = "'9"
$ws.Range("B20").Value = 5810538568
$ws.Range("C20").Value = "2021-11-22 06:45:54"
$ws.Range("D20").Value = "冬の伝言"
$ws.Range("E20").Value = "唉，怎么说呢，自从进击的轨迹之后近年的陛下仿佛是换了一种形象，虽说是放开了许多，但也让人感到缺失了5.6.7平表演时的悲伤、感动和纪行时的那种温柔、坚毅、富有人格魅力的感觉"
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 0

$ws.Range("A21").Value = "'9.1"
$ws.Range("B21").Value = 5811776235
$ws.Range("C21").Value = "2021-11-22 12:37:31"
$ws.Range("D21").Value = "不是你的朱雀"
$ws.Range("E21").Value = "从忧郁小王子变成了阳光大男孩"
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 5810538568

$ws.Range("A22").Value = "'9.2"
$ws.Range("B22").Value = 5814346853
$ws.Range("C22").Value = "2021-11-22 20:52:20"
$ws.Range("D22").Value = "Ponster_"
$ws.Range("E22").Value = "说起来陛下已经是40+的中年大叔了啊（大不敬）"
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 5810538568

$ws.Range("A23").Value = "'8"
$ws.Range("B23").Value = 5812089382
$ws.Range("C23").Value = "2021-11-22 13:36:01"
$ws.Range("D23").Value = "enemin"
$ws.Range("E23").Value = "太感谢了 之前还有一版安可不知道有没有大佬传 几乎是猫咪铃唱人偶的"
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 0

$ws.Range("A24").Value = "'8.1"
$ws.Range("B24").Value = 5814158584
$ws.Range("C24").Value = "2021-11-22 20:21:43"
$ws.Range("D24").Value = "林花花花"
$ws.Range("E24").Value = "那个b站之前有，是被删了吗"
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = 5812089382

$ws.Range("A25").Value = "'8.2"
$ws.Range("B25").Value = 5814154479
$ws.Range("C25").Value = "2021-11-22 20:22:06"
$ws.Range("D25").Value = "林花花花"
$ws.Range("E25").Value = "还好缓存的快"
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 5812089382

$ws.Range("A26").Value = "'8.3"
$ws.Range("B26").Value = 5814448893
$ws.Range("C26").Value = "2021-11-22 21:08:40"
$ws.Range("D26").Value = "enemin"
$ws.Range("E26").Value = "回复 @林花花花 :是啊 我前一秒还在看 然后推出去发现就无了 没有缓存 伤心了"
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 0
$ws.Range("H26").Value = 5812089382

$ws.Range("A27").Value = "'7"
$ws.Range("B27").Value = 5812890404
$ws.Range("C27").Value = "2021-11-22 16:48:22"
$ws.Range("D27").Value = "大白梨°"
$ws.Range("E27").Value = "[doge]陛下终究还是老了，不装13我很不习惯的"
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 0

$ws.Range("A28").Value = "'6"
$ws.Range("B28").Value = 5813229158
$ws.Range("C28").Value = "2021-11-22 17:57:18"
$ws.Range("D28").Value = "林花花花"
$ws.Range("E28").Value = "感谢！！"
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 0

$ws.Range("A29").Value = "'6.1"
$ws.Range("B29").Value = 5813402664
$ws.Range("C29").Value = "2021-11-22 18:25:22"
$ws.Range("D29").Value = "Ponster_"
$ws.Range("E29").Value = "也感谢你(=・ω・=)"
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 5813229158

$ws.Range("A30").Value = "'6.2"
$ws.Range("B30").Value = 5814183617
$ws.Range("C30").Value = "2021-11-22 20:26:08"
$ws.Range("D30").Value = "林花花花"
$ws.Range("E30").Value = "回复 @Ponster_ :想看八平[笑哭]"
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 5813229158

$ws.Range("A31").Value = "'5"
$ws.Range("B31").Value = 5815080989
$ws.Range("C31").Value = "2021-11-22 22:43:19"
$ws.Range("D31").Value = "艾奥萝卜"
$ws.Range("E31").Value = "先马再看"
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 0

$ws.Range("A32").Value = "'4"
$ws.Range("B32").Value = 5815188052
$ws.Range("C32").Value = "2021-11-22 22:59:05"
$ws.Range("D32").Value = "自律-Official"
$ws.Range("E32").Value = "哦哦哦赶紧缓存爽到[夏诺雅_太会了]"
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 2
$ws.Range("H32").Value = 0

$ws.Range("A33").Value = "'4.1"
$ws.Range("B33").Value = 5815284440
$ws.Range("C33").Value = "2021-11-22 23:12:58"
$ws.Range("D33").Value = "Ponster_"
$ws.Range("E33").Value = "这里不要脸地对简介链接里的仓库进行一个安利[doge]`n视频的评论区、弹幕已备份到其中，刚刚更新过`n[吃瓜][吃瓜][吃瓜]"
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 0
$ws.Range("H33").Value = 5815188052

$ws.Range("A34").Value = "'4.2"
$ws.Range("B34").Value = 5815389687
$ws.Range("C34").Value = "2021-11-22 23:28:15"
$ws.Range("D34").Value = "自律-Official"
$ws.Range("E34").Value = "回复 @Ponster_ :好起来了[夏诺雅_震撼]"
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 5815188052

$ws.Range("A35").Value = "'3"
$ws.Range("B35").Value = 5815263944
$ws.Range("C35").Value = "2021-11-22 23:09:47"
$ws.Range("D35").Value = "敦肃皇贵妃葛小队"
$ws.Range("E35").Value = "第一次看西装男吓尿了"
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 1
$ws.Range("H35").Value = 0

$ws.Range("A36").Value = "'3.1"
$ws.Range("B36").Value = 5815293406
$ws.Range("C36").Value = "2021-11-22 23:14:02"
$ws.Range("D36").Value = "Ponster_"
$ws.Range("E36").Value = "陛下可是“音乐界的杀人贵公子”呢"
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 5815263944

$ws.Range("A37").Value = "'20"
$ws.Range("B37").Value = 5807507484
$ws.Range("C37").Value = "2021-11-21 18:22:28"
$ws.Range("D37").Value = "syyuansang"
$ws.Range("E37").Value = "这是第几场，泪目了"
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 0

$ws.Range("B1").Copy()
$ws.Range("A2:A37").PasteSpecial(-4122)
$excel.CutCopyMode = 0